# Applies the "make setup was run" edit:
#  1. Reworks the "Dataset" sheet headers from
#       primary_email, birth_date, age_in_years, vital_status, id, name, description
#     down to
#       id, name, description
#     and removes the vital_status list data-validation.
#  2. Appends four new (empty) sheets after "DatasetCollection":
#       Sample, Molecule, Chemical_Substance, Analysis

$wb = $excel.ActiveWorkbook

# --- 1. Rebuild the "Dataset" sheet header row -----------------------------
$dataset = $wb.Worksheets.Item("Dataset")

# Drop any existing data validation rule(s) tied to this sheet (the
# vital_status ALIVE/DEAD/UNKNOWN dropdown lived on column D).
$dataset.Cells.Validation.Delete()

# Clear all existing cell content (old 7-column header row) and write the
# new 3-column header row.
$dataset.Cells.Clear()
$dataset.Range("A1").Value = "id"
$dataset.Range("B1").Value = "name"
$dataset.Range("C1").Value = "description"

# --- 2. Append the four new sheets ------------------------------------------
$newSheetNames = @("Sample", "Molecule", "Chemical_Substance", "Analysis")

foreach ($sheetName in $newSheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet.Name = $sheetName
}
